# BOT; UPDATE DATA
# Insert a new daily-data row (2020-05-05) above the trailing footnote row on
# the "相談件数" sheet, pushing the footnote from row 101 to row 102, then
# refresh the sheet's Print_Area defined name to cover the extra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row at 101 (shifts the footnote row 101 -> 102, and copies
# the number formatting of the row above, matching rows 99/100).
$ws.Rows("101:101").Insert()

# Fill in the new day's figures.
$ws.Range("A101").Value = 43956
$ws.Range("B101").Value = 321
$ws.Range("C101").Value = 33785
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 6958

# Extend the printable range to include the newly inserted row.
$wb.Names.Item("相談件数!Print_Area").RefersTo = '=相談件数!$A$1:$E$103'

# Move the active selection to the new last data cell.
[void]$ws.Range("B102").Select()
